$d = $word.ActiveDocument

# 1. Split the "So I got to work..." paragraph into two runs: append a
#    trailing space to the existing sentence and add a new sentence about
#    the Consolas font change as a new run in the same paragraph.
$oldText = "So I got to work and tried to create a theme that was black and green to simulate a bit of a terminal look... unfortunately for me, the outlines added how my code was handling it was completely off and it was pretty rough. On the other hand, I was able to find a style I wanted after changing twice which was still green and black but with a bit more of a different less complex look."
$newText = "So I got to work and tried to create a theme that was black and green to simulate a bit of a terminal look... unfortunately for me, the outlines added how my code was handling it was completely off and it was pretty rough. On the other hand, I was able to find a style I wanted after changing twice which was still green and black but with a bit more of a different less complex look. I also changed the font to Consolas which is the same font that VS Code uses to give a sense of familiarity."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newText, 2)

Write-Output "done"
